$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be stored as text, matching the source data
# (values like "29.748.27" / "1.000" / "0.9999" are text, not numbers)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.721.41"
$ws.Range("E2").Value = "  +8.34%  "
$ws.Range("D3").Value = "1.947.63"
$ws.Range("E3").Value = "  +6.88%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "342.17"
$ws.Range("E5").Value = "  +2.86%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("E7").Value = "  +4.37%  "
$ws.Range("D8").Value = "0.4144"
$ws.Range("E8").Value = "  +8.79%  "
$ws.Range("D9").Value = "48.37"
$ws.Range("E9").Value = "  +5.20%  "
$ws.Range("D10").Value = "0.08258"
$ws.Range("D11").Value = "1.042"
$ws.Range("E11").Value = "  +8.85%  "
$ws.Range("D12").Value = "22.65"
$ws.Range("E12").Value = "  +7.77%  "
$ws.Range("D13").Value = "1.937.24"
$ws.Range("E13").Value = "  +5.78%  "
$ws.Range("D14").Value = "6.198"
$ws.Range("E14").Value = "  +6.31%  "
$ws.Range("D15").Value = "7.425"
$ws.Range("E15").Value = "  +5.00%  "
$ws.Range("D16").Value = "92.31"
$ws.Range("E16").Value = "  +3.36%  "
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("E18").Value = "  +4.51%  "
$ws.Range("D19").Value = "0.06665"
$ws.Range("E19").Value = "  +1.27%  "
$ws.Range("E20").Value = "  +5.83%  "
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").Value = "29.695.24"
$ws.Range("E22").Value = "  +8.30%  "
$ws.Range("D23").Value = "5.604"
$ws.Range("E23").Value = "  +5.98%  "
$ws.Range("D24").Value = "11.25"
$ws.Range("E24").Value = "  +4.02%  "
$ws.Range("D25").Value = "2.285"
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("D26").Value = "2.171.39"
$ws.Range("E26").Value = "  +6.01%  "
$ws.Range("D27").Value = "160.73"
$ws.Range("E27").Value = "  +2.22%  "
$ws.Range("D28").Value = "20.17"
$ws.Range("E28").Value = "  +4.36%  "
$ws.Range("D29").Value = "2.195"
$ws.Range("E29").Value = "  +7.58%  "
$ws.Range("D30").Value = "5.680"
$ws.Range("E30").Value = "  +8.04%  "
$ws.Range("D31").Value = "122.45"
$ws.Range("E31").Value = "  +4.04%  "
$ws.Range("D32").Value = "1.024"
$ws.Range("E32").Value = "  +10.13%  "
$ws.Range("D33").Value = "0.09641"
$ws.Range("E33").Value = "  +3.63%  "
$ws.Range("D34").Value = "1.477"
$ws.Range("E34").Value = "  +12.46%  "
$ws.Range("D35").Value = "3.682"
$ws.Range("E35").Value = "  +3.34%  "
$ws.Range("D36").Value = "5.497"
$ws.Range("E36").Value = "  +5.45%  "
$ws.Range("E37").Value = "  +6.82%  "
$ws.Range("D38").Value = "0.02328"
$ws.Range("E38").Value = "  +6.49%  "
$ws.Range("D39").Value = "8.589"
$ws.Range("E39").Value = "  +6.08%  "
$ws.Range("D40").Value = "1.200"
$ws.Range("E40").Value = "  +5.29%  "
$ws.Range("D41").Value = "0.6120"
$ws.Range("E41").Value = "  +6.66%  "
$ws.Range("E42").Value = "  +8.28%  "
$ws.Range("E43").Value = "  +4.94%  "
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "1.272"
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "2.386"
$ws.Range("E46").Value = "  +32.78%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "0.5716"
$ws.Range("E47").Value = "  +6.08%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "12.52"
$ws.Range("E48").Value = "  +6.46%  "
$ws.Range("D49").Value = "2.004"
$ws.Range("E49").Value = "  +7.03%  "
$ws.Range("D50").Value = "0.07390"
$ws.Range("E50").Value = "  +12.45%  "
$ws.Range("D51").Value = "114.18"
$ws.Range("E51").Value = "  +3.57%  "
